$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-11, columns A, B, H, I, J, K.
# (Row 3 is untouched by this edit, so it is omitted below.)
# Numeric-looking text values (company numbers in B, SIC codes in I) are
# written with a leading apostrophe so the engine stores them as TEXT
# (matching the source file's inlineStr/text cells) instead of coercing
# them to numbers. The Style reset afterwards clears the transient
# "quote prefix" formatting so no stray style is left on the cell.

$rows = @(
    @{ Row = 2;  A = "GANDER INVESTMENTS LTD"; B = "16473515"; H = "Investments"; I = "68100,68209"; J = ""; K = "" },
    @{ Row = 4;  A = "SEVEN (HOLDCO) LIMITED"; B = "16473606"; H = "Other"; I = "64209"; J = "Activities of other holding companies n.e.c."; K = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles." },
    @{ Row = 5;  A = "INTERCONTINENTAL HOLDING COMPANY LIMITED"; B = "16473418"; H = "Other"; I = "64209"; J = "Activities of other holding companies n.e.c."; K = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles." },
    @{ Row = 6;  A = "AJ INVESTMENT AND CONSULTANCY LTD"; B = "16473328"; H = "Investments"; I = "64306,70229"; J = "Activities of real estate investment trusts"; K = "UK-regulated REIT companies." },
    @{ Row = 7;  A = "GAUNT CAPITAL LTD"; B = "16473262"; H = "Capital"; I = "64209"; J = "Activities of other holding companies n.e.c."; K = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles." },
    @{ Row = 8;  A = "TLJ INVESTMENT LTD"; B = "16473151"; H = "Investments"; I = "41100,55100,68100"; J = ""; K = "" },
    @{ Row = 9;  A = "THE DISLEY GROUP LTD"; B = "16473398"; H = "Other"; I = "64209"; J = "Activities of other holding companies n.e.c."; K = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles." },
    @{ Row = 10; A = "MARMIMI HOLDING LIMITED"; B = "16473234"; H = "Other"; I = "64209"; J = "Activities of other holding companies n.e.c."; K = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles." },
    @{ Row = 11; A = "BRIDGEWICK PARTNERS LIMITED"; B = "16473142"; H = "Partners"; I = "64999"; J = "Financial intermediation not elsewhere classified"; K = "Catch-all credit-oriented SPVs for novel lending structures." }
)

foreach ($item in $rows) {
    $r = $item.Row

    # Column A - company name (plain text, never numeric-looking)
    $ws.Cells.Item($r, 1).Value = $item.A

    # Column B - company number (numeric-looking -> force text)
    $ws.Cells.Item($r, 2).Value = "'" + $item.B
    $ws.Cells.Item($r, 2).Style = "Normal"

    # Column H - category (plain text)
    $ws.Cells.Item($r, 8).Value = $item.H

    # Column I - SIC codes (numeric-looking -> force text)
    $ws.Cells.Item($r, 9).Value = "'" + $item.I
    $ws.Cells.Item($r, 9).Style = "Normal"

    # Column J - SIC description (plain text, possibly empty)
    if ($item.J -eq "") {
        $ws.Cells.Item($r, 10).Value = "'"
        $ws.Cells.Item($r, 10).Style = "Normal"
    } else {
        $ws.Cells.Item($r, 10).Value = $item.J
    }

    # Column K - typical use case (plain text, possibly empty)
    if ($item.K -eq "") {
        $ws.Cells.Item($r, 11).Value = "'"
        $ws.Cells.Item($r, 11).Style = "Normal"
    } else {
        $ws.Cells.Item($r, 11).Value = $item.K
    }
}
